$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '28.131.92'
$ws.Range('E2').Value = '  -3.19%  '
$ws.Range('D3').Value = '1.912.15'
$ws.Range('E3').Value = '  -3.95%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -1.13%  '
$ws.Range('D5').Value = "'328.50"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = "'1.003"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('E7').Value = '  -6.27%  '
$ws.Range('D8').Value = "'0.4007"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.08%  '
$ws.Range('D9').Value = "'53.14"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.79%  '
$ws.Range('D10').Value = "'0.08390"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.31%  '
$ws.Range('E11').Value = '  -4.12%  '
$ws.Range('D12').Value = "'22.05"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.48%  '
$ws.Range('D13').Value = '1.889.23'
$ws.Range('E13').Value = '  -4.88%  '
$ws.Range('D14').Value = "'7.408"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.78%  '
$ws.Range('D15').Value = "'6.050"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.29%  '
$ws.Range('D16').Value = "'1.004"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').Value = "'89.45"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('D18').Value = "'0.00001065"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.19%  '
$ws.Range('D19').Value = "'0.06604"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('D20').Value = "'17.86"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'1.003"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').Value = "'5.740"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.56%  '
$ws.Range('D23').Value = '28.131.92'
$ws.Range('E24').Value = '  -6.48%  '
$ws.Range('D25').Value = "'2.303"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('D26').Value = '2.140.71'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').Value = "'152.96"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('D28').Value = "'20.02"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.62%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = "'5.737"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.82%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = "'2.124"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.90%  '
$ws.Range('D31').Value = "'123.36"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = "'0.09640"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'0.9734"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.33%  '
$ws.Range('E34').Value = '  -5.59%  '
$ws.Range('D35').Value = "'3.645"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').Value = "'5.540"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.95%  '
$ws.Range('D37').Value = "'1.274"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.63%  '
$ws.Range('D38').Value = "'8.802"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('D39').Value = "'0.02295"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.63%  '
$ws.Range('D40').Value = "'0.06148"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.44%  '
$ws.Range('D41').Value = "'0.6156"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('D42').Value = "'10.97"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.52%  '
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('D44').Value = "'0.1905"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.44%  '
$ws.Range('D45').Value = "'1.311"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5856"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.90%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'12.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('E48').Value = '  -6.41%  '
$ws.Range('D49').Value = "'3.437"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').Value = "'0.06916"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('D51').Value = "'111.44"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.30%  '
